$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.186.85"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "1.859.04"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.623"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.53%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.86"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.59%  "
$ws.Range("E9").Value = "  +2.14%  "
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("D12").Value = "2.127.73"
$ws.Range("E12").Value = "  +2.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "1.860.69"
$ws.Range("E14").Value = "  +2.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.677"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("D17").Value = "35.170.14"
$ws.Range("E17").Value = "  +1.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.19%  "
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +28.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.59%  "
$ws.Range("E28").Value = "  +2.25%  "
$ws.Range("E29").Value = "  +0.28%  "
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0557"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.83%  "
$ws.Range("E32").Value = "  +2.93%  "
$ws.Range("E33").Value = "  +28.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.834"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +20.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.60%  "
$ws.Range("E37").Value = "  +6.55%  "
$ws.Range("E38").Value = "  +7.60%  "
$ws.Range("E39").Value = "  +4.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "89.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.70%  "
$ws.Range("D41").Value = "1.341.07"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.21%  "
$ws.Range("E43").Value = "  +4.06%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0573"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.93%  "
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("B47").Value = "Gas"
$ws.Range("C47").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +43.38%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.28%  "
$ws.Range("D49").Value = "2.041.83"
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("E50").Value = "  +1.88%  "
$ws.Range("E51").Value = "  +0.31%  "
